$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sourcefield")

# Insert two new columns at the front (A and B), shifting existing columns right
$ws.Range("A1:B1").EntireColumn.Insert()

$ws.Range("A1").Value = "Source_project_id"
$ws.Range("B1").Value = "Source_dataset"
$ws.Range("A1:B1").Interior.Color = $ws.Range("C1").Interior.Color

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = "i-ier1-6j336sl3-h9urmye1jqo7ms"
    $ws.Cells.Item($r, 2).Value = "dbt_lend"
}

$ws.Range("A2:A10").WrapText = $true
$ws.Range("A2:B10").RowHeight = 29.25

$ws.Columns.Item(1).ColumnWidth = 27.6
$ws.Columns.Item(2).ColumnWidth = 20

$ws.Range("A2:B10").Select()
